$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1361.3334
$ws.Range("J19").Value = 897
$ws.Range("L19").Value = 897
$ws.Range("N19").Value = -1247

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 50003376
$ws.Range("I64").Value = 50003376
$ws.Range("K64").Value = 50003376
$ws.Range("M64").Value = -50003128

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 50003376
$ws.Range("I67").Value = 50003376
$ws.Range("K67").Value = 50003376
$ws.Range("M67").Value = -50002518

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 20879208
$ws.Range("I74").Value = 20879208
$ws.Range("K74").Value = 20879208
$ws.Range("M74").Value = -20878272

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 20879208
$ws.Range("I77").Value = 20879208
$ws.Range("K77").Value = 104396040
$ws.Range("M77").Value = -104391360

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 9353.5
$ws.Range("I86").Value = 9204.695
$ws.Range("J86").Value = 9664.637000000001
$ws.Range("K86").Value = 9204.695
$ws.Range("L86").Value = 9664.637000000001
$ws.Range("M86").Value = -8081.695
$ws.Range("N86").Value = -11910.637

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 9353.5
$ws.Range("I89").Value = 9204.695
$ws.Range("J89").Value = 9664.637000000001
$ws.Range("K89").Value = 46023.475
$ws.Range("L89").Value = 48323.185
$ws.Range("M89").Value = -40407.475
$ws.Range("N89").Value = -59555.185

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2571.3333
$ws.Range("J97").Value = 2571.3333
$ws.Range("L97").Value = 7713.999899999999
$ws.Range("N97").Value = -8705.999899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3229.5
$ws.Range("I137").Value = 2939.4614
$ws.Range("K137").Value = 8818.3842
$ws.Range("M137").Value = -6268.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1866
$ws.Range("I45").Value = 1649.5
$ws.Range("K45").Value = 1649.5
$ws.Range("M45").Value = -1272.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 20000
$ws.Range("J19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("N19").Value = -20346

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 100002360
$ws.Range("I134").Value = 100002360
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 300007080
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = -300004545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4035.111
$ws.Range("I86").Value = 3865.2
$ws.Range("K86").Value = 3865.2
$ws.Range("M86").Value = -2742.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4035.111
$ws.Range("I89").Value = 3865.2
$ws.Range("K89").Value = 19326
$ws.Range("M89").Value = -13710

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 932.3333
$ws.Range("I94").Value = 913.4286
$ws.Range("K94").Value = 913.4286
$ws.Range("M94").Value = -462.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2225218
$ws.Range("I105").Value = 2858994.5
$ws.Range("K105").Value = 2858994.5
$ws.Range("M105").Value = -2857247.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 15626110
$ws.Range("I134").Value = 20834238
$ws.Range("K134").Value = 62502714
$ws.Range("M134").Value = -62500179

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 388984
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 913.3333
$ws.Range("I34").Value = 696
$ws.Range("K34").Value = 2088
$ws.Range("M34").Value = -2004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3333
$ws.Range("J39").Value = 8000
$ws.Range("L39").Value = 24000
$ws.Range("N39").Value = -24588

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2147.25
$ws.Range("I63").Value = 1996.3334
$ws.Range("J63").Value = 2600
$ws.Range("K63").Value = 5989.0002
$ws.Range("L63").Value = 7800
$ws.Range("M63").Value = -5240.0002
$ws.Range("N63").Value = -9298

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 2147.25
$ws.Range("I66").Value = 1996.3334
$ws.Range("J66").Value = 2600
$ws.Range("K66").Value = 17967.0006
$ws.Range("L66").Value = 23400
$ws.Range("M66").Value = -14223.0006
$ws.Range("N66").Value = -30888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 12882.083
$ws.Range("I70").Value = 5765.8335
$ws.Range("K70").Value = 17297.5005
$ws.Range("M70").Value = -16982.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 12882.083
$ws.Range("I73").Value = 5765.8335
$ws.Range("K73").Value = 17297.5005
$ws.Range("M73").Value = -16205.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1598.25
$ws.Range("I70").Value = 1598.25
$ws.Range("K70").Value = 1598.25
$ws.Range("M70").Value = -1328.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 1598.25
$ws.Range("I73").Value = 1598.25
$ws.Range("K73").Value = 1598.25
$ws.Range("M73").Value = -662.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 978.9048
$ws.Range("I97").Value = 626.36365
$ws.Range("K97").Value = 626.36365
$ws.Range("M97").Value = -130.36365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 740.6667
$ws.Range("I102").Value = 740.6667
$ws.Range("K102").Value = 740.6667
$ws.Range("M102").Value = 881.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 623.26666
$ws.Range("I107").Value = 282
$ws.Range("J107").Value = 1305.8
$ws.Range("K107").Value = 282
$ws.Range("L107").Value = 1305.8
$ws.Range("M107").Value = 1638
$ws.Range("N107").Value = -5145.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 61784.293
$ws.Range("I113").Value = 68822.266
$ws.Range("J113").Value = 8999.5
$ws.Range("K113").Value = 68822.266
$ws.Range("L113").Value = 8999.5
$ws.Range("M113").Value = -66652.266
$ws.Range("N113").Value = -13339.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5112.375
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 6166.5
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 18499.5
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -23399.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2841.8
$ws.Range("I16").Value = 717.25
$ws.Range("J16").Value = 4258.1665
$ws.Range("K16").Value = 717.25
$ws.Range("L16").Value = 4258.1665
$ws.Range("M16").Value = -547.25
$ws.Range("N16").Value = -4598.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3879.7
$ws.Range("I46").Value = 1466.1666
$ws.Range("K46").Value = 1466.1666
$ws.Range("M46").Value = -1278.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 232.33333
$ws.Range("I55").Value = 210
$ws.Range("K55").Value = 210
$ws.Range("M55").Value = -37

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1321.0625
$ws.Range("I93").Value = 1321.0625
$ws.Range("K93").Value = 1321.0625
$ws.Range("M93").Value = -73.0625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 17000
$ws.Range("J106").Value = 17000
$ws.Range("L106").Value = 17000
$ws.Range("N106").Value = -19524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11642.4
$ws.Range("I122").Value = 11642.4
$ws.Range("K122").Value = 34927.2
$ws.Range("M122").Value = -32477.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2261.5
$ws.Range("I81").Value = 2451.3333
$ws.Range("J81").Value = 1122.5
$ws.Range("K81").Value = 4902.6666
$ws.Range("L81").Value = 2245
$ws.Range("M81").Value = -3841.6666
$ws.Range("N81").Value = -4367

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2261.5
$ws.Range("I84").Value = 2451.3333
$ws.Range("J84").Value = 1122.5
$ws.Range("K84").Value = 24513.333
$ws.Range("L84").Value = 11225
$ws.Range("M84").Value = -19209.333
$ws.Range("N84").Value = -21833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 33998.668
$ws.Range("J97").Value = 33998.668
$ws.Range("L97").Value = 33998.668
$ws.Range("N97").Value = -35980.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 82125.664
$ws.Range("J106").Value = 82125.664
$ws.Range("L106").Value = 82125.664
$ws.Range("N106").Value = -84649.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13162937
$ws.Range("I132").Value = 20002386
$ws.Range("J132").Value = 10148.923
$ws.Range("K132").Value = 60007158
$ws.Range("L132").Value = 30446.769
$ws.Range("M132").Value = -60004628
$ws.Range("N132").Value = -35506.769
